$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fjernet mellomrom" (removed whitespace): trailing space removed from
# the value in E7 ("Teltplass lille Foretjørn ") -> ("Teltplass lille Foretjørn")
$ws.Range("E7").Value = "Teltplass lille Foretjørn"
